# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns
# with the latest figures from the scheduled GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price / Volume(1h) values per row, keyed by cell address.
$updates = @{
    "D2" = "254.84"
    "E2" = "3.52%"
    "D3" = "28.20"
    "E3" = "-5.23%"
    "D4" = "5.327"
    "E4" = "3.19%"
    "D5" = "0.05847"
    "E5" = "0.76%"
    "D6" = "6.709"
    "E6" = "0.41%"
    "D7" = "0.8663"
    "E7" = "1.68%"
    "D8" = "0.9118"
    "E8" = "5.76%"
    "D9" = "0.1422"
    "E9" = "3.30%"
    "D10" = "0.07165"
    "E10" = "0.87%"
    "D11" = "0.03180"
    "E11" = "-0.78%"
    "D12" = "0.09219"
    "E12" = "-1.69%"
    "D13" = "0.001537"
    "E13" = "0.03%"
    "D14" = "0.0006057"
    "E14" = "1.28%"
    "D15" = "0.005817"
    "E15" = "-1.26%"
    "D16" = "3.499"
    "E16" = "0.03%"
    "D17" = "3.231"
    "E17" = "0.11%"
    "D18" = "2.201"
    "E18" = "-0.44%"
    "D19" = "0.3171"
    "E19" = "-0.80%"
    "D20" = "0.03447"
    "E20" = "2.99%"
    "D21" = "0.1316"
    "E21" = "1.17%"
    "D22" = "3.541"
    "E22" = "1.43%"
    "D23" = "0.04166"
    "E23" = "0.52%"
    "E24" = "-0.16%"
    "D25" = "0.005042"
    "E25" = "21.69%"
    "D26" = "0.001228"
    "E26" = "0.07%"
    "D27" = "0.0001199"
    "E27" = "-0.06%"
    "D28" = "0.0001937"
    "E28" = "34.02%"
    "D40" = "0.03852"
    "E40" = "2.43%"
    "D41" = "0.1101"
    "E41" = "2.86%"
    "D42" = "0.002388"
    "E42" = "8.57%"
    "D43" = "0.002947"
    "E43" = "-48.09%"
    "D44" = "0.01092"
    "E44" = "14.32%"
    "D45" = "0.00005230"
    "E45" = "-1.25%"
    "E46" = "-0.02%"
    "D47" = "0.08981"
    "E47" = "54.83%"
    "D48" = "0.002155"
    "E48" = "-1.21%"
    "D49" = "0.00002099"
    "E49" = "-0.02%"
    "D50" = "0.0001999"
    "E50" = "-0.02%"
}

foreach ($addr in $updates.Keys) {
    # Leading apostrophe forces text (matching the source sheet, where
    # Price/Volume are stored as strings, e.g. "3.52%"), not a number.
    $ws.Range($addr).Value = "'" + $updates[$addr]
    # Drop back to the Normal style so the text-entry quote-prefix
    # formatting doesn't leave a stray style on the cell.
    $ws.Range($addr).Style = "Normal"
}
